$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Swap the match-detail columns (F:V) between the following row
# pairs. Columns A:E (Indice, pais, torneio, temporada, data_partida)
# stay attached to their original row.
# ---------------------------------------------------------------
function Swap-Rows($r1, $r2) {
    $cols = 6..22   # F=6 .. V=22
    foreach ($c in $cols) {
        $v1 = $ws.Cells.Item($r1, $c).Value2
        $v2 = $ws.Cells.Item($r2, $c).Value2
        $ws.Cells.Item($r1, $c).Value2 = $v2
        $ws.Cells.Item($r2, $c).Value2 = $v1
    }
}

Swap-Rows 39 40
Swap-Rows 43 44
Swap-Rows 53 54

# ---------------------------------------------------------------
# Append new row 68 (new betexplorer match record), copying the
# cell formatting from row 67 and then filling in the values.
# ---------------------------------------------------------------
$ws.Range("A67:V67").Copy()
$ws.Range("A68:V68").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(68, 1).Value2 = 67
$ws.Cells.Item(68, 2).Value2 = "spain"
$ws.Cells.Item(68, 3).Value2 = "laliga2"
$ws.Cells.Item(68, 4).Value2 = "2023-2024"
$ws.Cells.Item(68, 5).Value2 = 45191.875
$ws.Cells.Item(68, 6).Value2 = "FC Cartagena SAD"
$ws.Cells.Item(68, 7).Value2 = 1
$ws.Cells.Item(68, 8).Value2 = "Eibar"
$ws.Cells.Item(68, 9).Value2 = 2
$ws.Cells.Item(68, 10).Value2 = 2.86
$ws.Cells.Item(68, 11).Value2 = "17/09/2023 15:42"
$ws.Cells.Item(68, 12).Value2 = 3.51
$ws.Cells.Item(68, 13).Value2 = "22/09/2023 20:32"
$ws.Cells.Item(68, 14).Value2 = 3.01
$ws.Cells.Item(68, 15).Value2 = "17/09/2023 15:42"
$ws.Cells.Item(68, 16).Value2 = 3.08
$ws.Cells.Item(68, 17).Value2 = "22/09/2023 20:32"
$ws.Cells.Item(68, 18).Value2 = 2.86
$ws.Cells.Item(68, 19).Value2 = "17/09/2023 15:42"
$ws.Cells.Item(68, 20).Value2 = 2.39
$ws.Cells.Item(68, 21).Value2 = "22/09/2023 20:32"
$ws.Cells.Item(68, 22).Value2 = "https://www.betexplorer.com/football/spain/laliga2/fc-cartagena-sad-eibar/f7VXoLKe/"
